$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4. This shifts the existing rows 4..71 down
# to 5..72 (matching the diff, where every existing record after row 3
# moves down by one row).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record from the diff.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(4, 3).Value = "Coquimbo"
$ws.Cells.Item(4, 4).Value = 44882
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 100112022
$ws.Cells.Item(4, 7).Value = "Arveja Verde"
$ws.Cells.Item(4, 8).Value = "Perfection"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 65
$ws.Cells.Item(4, 11).Value = 19000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 19462
$ws.Cells.Item(4, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 778
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
